# balance: Lower zombie melee skill. Update some descriptions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stuff-descriptor")

# New description for Player (row 4, D4 was empty)
$ws.Range("D4").Value = "Your very self"

# LeatherArmor (row 9): replace computed description formula with static text
$ws.Range("D9").Value = "Comfy armor"

# Sword (row 11): replace computed description formula with static text
$ws.Range("D11").Value = "Larger weapon"

# Dagger (row 12): replace computed description formula with static text
$ws.Range("D12").Value = "Small weapon"

# RareSword (row 22): replace computed description formula with static text
$ws.Range("D22").Value = "Rare sword formally owned by a knight"

# RareDagger (row 23): replace computed description formula with static text
$ws.Range("D23").Value = "Rare dagger"

# Zombie (row 24): lower melee_skill from 3 to 2
$ws.Range("I24").Value = 2

# Update the active selection on the sheet
$ws.Range("I25").Select()
